$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 14 (the most recent "Ramas de apio" record) down to row 15
# so the existing record is preserved, then overwrite row 14 with this
# week's new reading (weekly fruit/vegetable price update).
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(15).Insert(-4121)

# Write the new weekly values into row 14
$ws.Range("D14").Value = 44504
$ws.Range("J14").Value = 55
$ws.Range("K14").Value = 4000
$ws.Range("L14").Value = 4000
$ws.Range("M14").Value = 4000
$ws.Range("P14").Value = 4000
